$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.277.50"
$ws.Range("E2").Value = "  +2.29%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.351.82"
$ws.Range("E3").Value = "  +6.07%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.76%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.74%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.642"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.00%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.12%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.637"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.77%  "

# Row 10 - Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.43%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +3.09%  "

# Row 12 - Polkadot
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.44%  "

# Row 13 - Polygon
$ws.Range("E13").Value = "  +7.22%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +2.14%  "

# Row 15 - Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.53%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "2.709.25"
$ws.Range("E16").Value = "  +6.34%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.434.82"
$ws.Range("E17").Value = "  +9.43%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "43.268.27"
$ws.Range("E18").Value = "  +2.30%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  +3.33%  "

# Row 20 - Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.03%  "

# Row 21 - Litecoin
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.16%  "

# Row 22 - ImmutableX (was PancakeSwap)
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.08%  "

# Row 23 - PancakeSwap (was ImmutableX)
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.88%  "

# Row 24 - BitcoinCash
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "254.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.54%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.13%  "

# Row 26 - Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.80%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.00%  "

# Row 28 - InjectiveProtocol
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.71%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  +1.16%  "

# Row 30 - EthereumClassic
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.75%  "

# Row 31 - Monero
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.04%  "

# Row 32 - WEMIXToken
$ws.Range("E32").Value = "  -1.06%  "

# Row 33 - Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0928"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.79%  "

# Row 34 - Filecoin
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.38%  "

# Row 35 - Stellar
$ws.Range("E35").Value = "  +5.27%  "

# Row 36 - RenderToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.85%  "

# Row 37 - NEARProtocol
$ws.Range("E37").Value = "  -5.37%  "

# Row 38 - VeChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0376"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.67%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  +1.93%  "

# Row 40 - LidoDAOToken
$ws.Range("E40").Value = "  +11.53%  "

# Row 41 - MultiversX
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.39%  "

# Row 42 - ARBITRUM
$ws.Range("E42").Value = "  +14.39%  "

# Row 43 - Algorand
$ws.Range("E43").Value = "  +0.90%  "

# Row 44 - Celestia
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.57%  "

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = "  +0.15%  "

# Row 46 - THORChain
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.87%  "

# Row 47 - FraxShare
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.49%  "

# Row 48 - Aave
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "111.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.56%  "

# Row 49 - TrustWalletToken
$ws.Range("E49").Value = "  -0.59%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  +3.61%  "

# Row 51 - WOONetwork
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.464"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.59%  "
